$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet ALC
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("ALC")

$ws.Range("H62").Value = 1406.7727
$ws.Range("I62").Value = 1397.3158
$ws.Range("J62").Value = 1466.6666
$ws.Range("K62").Value = 1397.3158
$ws.Range("L62").Value = 1466.6666
$ws.Range("M62").Value = -773.3158000000001
$ws.Range("N62").Value = -2714.6666

$ws.Range("H65").Value = 1406.7727
$ws.Range("I65").Value = 1397.3158
$ws.Range("J65").Value = 1466.6666
$ws.Range("K65").Value = 6986.579000000001
$ws.Range("L65").Value = 7333.333000000001
$ws.Range("M65").Value = -3866.579000000001
$ws.Range("N65").Value = -13573.333

# ---------------------------------------------------------------------------
# Sheet ARM
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("ARM")

$ws.Range("H2").Value = 1541.9744
$ws.Range("I2").Value = 793.7917
$ws.Range("J2").Value = 2739.0667
$ws.Range("K2").Value = 793.7917
$ws.Range("L2").Value = 2739.0667
$ws.Range("M2").Value = -680.7917
$ws.Range("N2").Value = -2965.0667

$ws.Range("H32").Value = 413943.9
$ws.Range("I32").Value = 4834.0166
$ws.Range("J32").Value = 2868603.2
$ws.Range("K32").Value = 4834.0166
$ws.Range("L32").Value = 2868603.2
$ws.Range("M32").Value = -4547.0166
$ws.Range("N32").Value = -2869177.2

$ws.Range("H116").Value = 1541.9744
$ws.Range("I116").Value = 793.7917
$ws.Range("J116").Value = 2739.0667
$ws.Range("K116").Value = 793.7917
$ws.Range("L116").Value = 2739.0667
$ws.Range("M116").Value = 1500.2083
$ws.Range("N116").Value = -7327.066699999999

# Rows 121-135 and 137-141 (row 136 left untouched): clear stale H:N data
$r = $ws.Range("H121:N135,H137:N141")
foreach ($area in $r.Areas) {
    $area.ClearContents()
}

# ---------------------------------------------------------------------------
# Sheet BSM
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("BSM")

$ws.Range("H3").Value = 1541.9744
$ws.Range("I3").Value = 793.7917
$ws.Range("J3").Value = 2739.0667
$ws.Range("K3").Value = 793.7917
$ws.Range("L3").Value = 2739.0667
$ws.Range("M3").Value = -679.7917
$ws.Range("N3").Value = -2967.0667

# ---------------------------------------------------------------------------
# Sheet GSM
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("GSM")

$ws.Range("H102").Value = 1999.6111
$ws.Range("I102").Value = 1875.375
$ws.Range("K102").Value = 1875.375
$ws.Range("M102").Value = -253.375

# ---------------------------------------------------------------------------
# Sheet WVR
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("WVR")

$ws.Range("H96").Value = 2758.1082
$ws.Range("I96").Value = 2189.5833
$ws.Range("J96").Value = 3807.6924
$ws.Range("K96").Value = 2189.5833
$ws.Range("L96").Value = 3807.6924
$ws.Range("M96").Value = -816.5832999999998
$ws.Range("N96").Value = -6553.6924
